$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the batter's name
$ws.Name = "Yuzvendra Chahal"

# Insert a new column A (shifts existing A:L -> B:M) to hold the new "matchNo" field.
# The shifted cells keep their original value types/styles automatically.
$ws.Columns.Item(1).Insert()

# Insert a new row 2 (shifts existing row 2 -> row 3) to hold the new match record.
$ws.Rows.Item(2).Insert()

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# ---- New row 2: the "31st" match record ----
$ws.Range("A2").Value = "31st"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Yuzvendra Chahal"

# "states" is blank text (matches the empty-but-text cell used elsewhere in the sheet)
# and "runs"/"balls"/"fours"/"sixes"/"sr" are numeric-looking values that must stay
# text, like the rest of the sheet. A leading apostrophe forces a literal text value
# (without re-parsing it as a number); resetting the style back to Normal afterward
# drops the quote-prefix formatting flag again so no visible formatting changes.
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'2"
$ws.Range("F2").Value = "'6"
$ws.Range("G2").Value = "'0"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'33.33"
$ws.Range("D2:I2").Style = "Normal"

$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "September 20"
$ws.Range("M2").Value = "KKR won by 9 wickets (with 60 balls remaining)"

# ---- Row 3 (the original row, shifted down): just fill the new matchNo cell ----
# All other cells in this row already carried their original values/types through
# the column/row insert shifts above, so only A3 is new.
$ws.Range("A3").Value = "19th"
